$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the values that were in column K into column D for rows 3-6, 8-14,
# and clear out column K for those rows.

$ws.Range("D3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K3").Value = ""

$ws.Range("D4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K4").Value = ""

$ws.Range("D5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K5").Value = ""

$ws.Range("D6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("K6").Value = ""

$ws.Range("D8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K8").Value = ""

$ws.Range("D9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K9").Value = ""

$ws.Range("D10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K10").Value = ""

$ws.Range("D11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("K11").Value = ""

$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
$ws.Range("K12").Value = ""

$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("K13").Value = ""

$ws.Range("D14").Value = "2.300 TL - 9.500 TL"
$ws.Range("K14").Value = ""
